$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- Sheet ALC ---
# row 103
$ws1.Range("H103").Value = 192432.23
$ws1.Range("I103").Value = 227401.81
$ws1.Range("J103").Value = 99.5
$ws1.Range("K103").Value = 682205.4299999999
$ws1.Range("L103").Value = 298.5
$ws1.Range("M103").Value = -681619.4299999999
$ws1.Range("N103").Value = -1470.5

# row 121
$ws1.Range("H121").Value = 870.26086
$ws1.Range("J121").Value = 870.26086
$ws1.Range("L121").Value = 2610.78258
$ws1.Range("N121").Value = -6104.78258

# row 132
$ws1.Range("H132").Value = 30446.084
$ws1.Range("I132").Value = 36105.7
$ws1.Range("J132").Value = 2148
$ws1.Range("K132").Value = 108317.1
$ws1.Range("L132").Value = 6444
$ws1.Range("M132").Value = -105787.1
$ws1.Range("N132").Value = -11504

# row 135
$ws1.Range("H135").Value = 20840182
$ws1.Range("I135").Value = 1181.8572
$ws1.Range("K135").Value = 10636.7148
$ws1.Range("M135").Value = -8101.7148

# row 137
$ws1.Range("H137").Value = 1712.5
$ws1.Range("I137").Value = 1228.5714
$ws1.Range("K137").Value = 3685.7142
$ws1.Range("M137").Value = -1135.7142

# row 138
$ws1.Range("H138").Value = 2995.6086
$ws1.Range("I138").Value = 472.5
$ws1.Range("J138").Value = 3526.7896
$ws1.Range("K138").Value = 1417.5
$ws1.Range("L138").Value = 10580.3688
$ws1.Range("M138").Value = 3722.5
$ws1.Range("N138").Value = -20860.3688

# row 141
$ws1.Range("H141").Value = 3099.1667
$ws1.Range("I141").Value = 1245
$ws1.Range("J141").Value = 4026.25
$ws1.Range("K141").Value = 3735
$ws1.Range("L141").Value = 12078.75
$ws1.Range("M141").Value = 1445
$ws1.Range("N141").Value = -22438.75

# --- Sheet ARM ---
# row 32
$ws2.Range("H32").Value = 2411.75
$ws2.Range("I32").Value = 1637.9584
$ws2.Range("J32").Value = 7054.5
$ws2.Range("K32").Value = 1637.9584
$ws2.Range("L32").Value = 7054.5
$ws2.Range("M32").Value = -1350.9584
$ws2.Range("N32").Value = -7628.5

# row 88
$ws2.Range("H88").Value = 58348.5
$ws2.Range("I88").Value = 1520.8334
$ws2.Range("J88").Value = 86762.336
$ws2.Range("K88").Value = 1520.8334
$ws2.Range("L88").Value = 86762.336
$ws2.Range("M88").Value = -1114.8334
$ws2.Range("N88").Value = -87574.336

# row 91
$ws2.Range("H91").Value = 58348.5
$ws2.Range("I91").Value = 1520.8334
$ws2.Range("J91").Value = 86762.336
$ws2.Range("K91").Value = 1520.8334
$ws2.Range("L91").Value = 86762.336
$ws2.Range("M91").Value = -116.8334
$ws2.Range("N91").Value = -89570.336

# row 105
$ws2.Range("H105").Value = 33491.668
$ws2.Range("J105").Value = 33491.668
$ws2.Range("L105").Value = 33491.668
$ws2.Range("N105").Value = -40479.668

# row 112
$ws2.Range("H112").Value = 32667.666
$ws2.Range("J112").Value = 32667.666
$ws2.Range("L112").Value = 32667.666
$ws2.Range("N112").Value = -35621.666

# row 132
$ws2.Range("H132").Value = 17364.469
$ws2.Range("I132").Value = 1361.4762
$ws2.Range("K132").Value = 4084.4286
$ws2.Range("M132").Value = -1554.4286

# --- Sheet BSM ---
# row 86
$ws3.Range("H86").Value = 1466.5
$ws3.Range("I86").Value = 1224.1666
$ws3.Range("J86").Value = 1951.1666
$ws3.Range("K86").Value = 1224.1666
$ws3.Range("L86").Value = 1951.1666
$ws3.Range("M86").Value = -101.1666
$ws3.Range("N86").Value = -4197.1666

# row 89
$ws3.Range("H89").Value = 1466.5
$ws3.Range("I89").Value = 1224.1666
$ws3.Range("J89").Value = 1951.1666
$ws3.Range("K89").Value = 6120.833000000001
$ws3.Range("L89").Value = 9755.833000000001
$ws3.Range("M89").Value = -504.8330000000005
$ws3.Range("N89").Value = -20987.833

# --- Sheet CRP ---
# row 7
$ws4.Range("H7").Value = 63.166668
$ws4.Range("J7").Value = 70
$ws4.Range("L7").Value = 70
$ws4.Range("N7").Value = -296

# row 31
$ws4.Range("H31").Value = 8721.392
$ws4.Range("I31").Value = 10451.03
$ws4.Range("J31").Value = 4330.769
$ws4.Range("K31").Value = 10451.03
$ws4.Range("L31").Value = 4330.769
$ws4.Range("M31").Value = -10156.03
$ws4.Range("N31").Value = -4920.769

# row 34
$ws4.Range("H34").Value = 8721.392
$ws4.Range("I34").Value = 10451.03
$ws4.Range("J34").Value = 4330.769
$ws4.Range("K34").Value = 10451.03
$ws4.Range("L34").Value = 4330.769
$ws4.Range("M34").Value = -10249.03
$ws4.Range("N34").Value = -4734.769

# row 43
$ws4.Range("H43").Value = 34609.332
$ws4.Range("J43").Value = 34609.332
$ws4.Range("L43").Value = 34609.332
$ws4.Range("N43").Value = -34977.332

# row 88
$ws4.Range("H88").Value = 30000
$ws4.Range("J88").Value = 30000
$ws4.Range("L88").Value = 30000
$ws4.Range("N88").Value = -30812

# row 91
$ws4.Range("H91").Value = 30000
$ws4.Range("J91").Value = 30000
$ws4.Range("L91").Value = 30000
$ws4.Range("N91").Value = -32808

# row 101
$ws4.Range("H101").Value = 34609.332
$ws4.Range("J101").Value = 34609.332
$ws4.Range("L101").Value = 34609.332
$ws4.Range("N101").Value = -41099.332

# row 134
$ws4.Range("H134").Value = 951.8182
$ws4.Range("I134").Value = 961.1111
$ws4.Range("J134").Value = 910
$ws4.Range("K134").Value = 2883.3333
$ws4.Range("L134").Value = 2730
$ws4.Range("M134").Value = -348.3332999999998
$ws4.Range("N134").Value = -7800

# --- Sheet CUL ---
# row 59
$ws5.Range("H59").Value = 2500
$ws5.Range("I59").Value = 2000
$ws5.Range("K59").Value = 6000
$ws5.Range("M59").Value = -5460

# row 68
$ws5.Range("H68").Value = 3521.1892
$ws5.Range("I68").Value = 590.4286
$ws5.Range("J68").Value = 7367.8125
$ws5.Range("K68").Value = 1771.2858
$ws5.Range("L68").Value = 22103.4375
$ws5.Range("M68").Value = -960.2857999999999
$ws5.Range("N68").Value = -23725.4375

# row 71
$ws5.Range("H71").Value = 3521.1892
$ws5.Range("I71").Value = 590.4286
$ws5.Range("J71").Value = 7367.8125
$ws5.Range("K71").Value = 5313.8574
$ws5.Range("L71").Value = 66310.3125
$ws5.Range("M71").Value = -1257.8574
$ws5.Range("N71").Value = -74422.3125

# row 75
$ws5.Range("I75").Value = 1500
$ws5.Range("J75").Value = 0
$ws5.Range("K75").Value = 4500
$ws5.Range("L75").Value = 0
$ws5.Range("M75").Value = -3502
$ws5.Range("N75").ClearContents()

# row 78
$ws5.Range("I78").Value = 1500
$ws5.Range("J78").Value = 0
$ws5.Range("K78").Value = 13500
$ws5.Range("L78").Value = 0
$ws5.Range("M78").Value = -8508
$ws5.Range("N78").ClearContents()

# row 81
$ws5.Range("H81").Value = 4216.8
$ws5.Range("J81").Value = 5119.375
$ws5.Range("L81").Value = 15358.125
$ws5.Range("N81").Value = -17604.125

# row 84
$ws5.Range("H84").Value = 4216.8
$ws5.Range("J84").Value = 5119.375
$ws5.Range("L84").Value = 46074.375
$ws5.Range("N84").Value = -57306.375

# row 87
$ws5.Range("H87").Value = 11085
$ws5.Range("I87").Value = 636
$ws5.Range("K87").Value = 1908
$ws5.Range("M87").Value = -660

# row 90
$ws5.Range("H90").Value = 11085
$ws5.Range("I90").Value = 636
$ws5.Range("K90").Value = 5724
$ws5.Range("M90").Value = 516

# row 131
$ws5.Range("H131").Value = 797.0700000000001
$ws5.Range("J131").Value = 823.33685
$ws5.Range("L131").Value = 2470.01055
$ws5.Range("N131").Value = -12550.01055

# --- Sheet GSM ---
# row 101
$ws6.Range("H101").Value = 29124.75
$ws6.Range("J101").Value = 29124.75
$ws6.Range("L101").Value = 29124.75
$ws6.Range("N101").Value = -35614.75

# row 104
$ws6.Range("H104").Value = 0
$ws6.Range("J104").Value = 0
$ws6.Range("L104").Value = 0
$ws6.Range("N104").ClearContents()

# --- Sheet LTW ---
# row 22
$ws7.Range("H22").Value = 959.3
$ws7.Range("I22").Value = 773
$ws7.Range("J22").Value = 1083.5
$ws7.Range("K22").Value = 773
$ws7.Range("L22").Value = 1083.5
$ws7.Range("M22").Value = -478
$ws7.Range("N22").Value = -1673.5

# row 27
$ws7.Range("H27").Value = 959.3
$ws7.Range("I27").Value = 773
$ws7.Range("J27").Value = 1083.5
$ws7.Range("K27").Value = 773
$ws7.Range("L27").Value = 1083.5
$ws7.Range("M27").Value = -666
$ws7.Range("N27").Value = -1297.5

# row 106
$ws7.Range("H106").Value = 19000
$ws7.Range("J106").Value = 19000
$ws7.Range("L106").Value = 19000
$ws7.Range("N106").Value = -21524

# row 107
$ws7.Range("H107").Value = 1666.3334
$ws7.Range("I107").Value = 1666.3334
$ws7.Range("K107").Value = 1666.3334
$ws7.Range("M107").Value = 253.6666

# row 132
$ws7.Range("H132").Value = 448115.38
$ws7.Range("I132").Value = 635637.75
$ws7.Range("J132").Value = 2749.75
$ws7.Range("K132").Value = 1906913.25
$ws7.Range("L132").Value = 8249.25
$ws7.Range("M132").Value = -1904383.25
$ws7.Range("N132").Value = -13309.25

# row 136
$ws7.Range("H136").Value = 2258.2856
$ws7.Range("I136").Value = 1561.5834
$ws7.Range("J136").Value = 3187.2222
$ws7.Range("K136").Value = 4684.7502
$ws7.Range("L136").Value = 9561.6666
$ws7.Range("M136").Value = -2134.7502
$ws7.Range("N136").Value = -14661.6666

# --- Sheet WVR ---
# row 62
$ws8.Range("H62").Value = 5600
$ws8.Range("J62").Value = 5600
$ws8.Range("L62").Value = 5600
$ws8.Range("N62").Value = -6848

# row 65
$ws8.Range("H65").Value = 5600
$ws8.Range("J65").Value = 5600
$ws8.Range("L65").Value = 28000
$ws8.Range("N65").Value = -34240

# row 98
$ws8.Range("H98").Value = 45000
$ws8.Range("J98").Value = 45000
$ws8.Range("L98").Value = 45000
$ws8.Range("N98").Value = -50990

# row 104
$ws8.Range("H104").Value = 31342.5
$ws8.Range("J104").Value = 31342.5
$ws8.Range("L104").Value = 31342.5
$ws8.Range("N104").Value = -38330.5

# row 131
$ws8.Range("H131").Value = 26443.666
$ws8.Range("J131").Value = 26443.666
$ws8.Range("L131").Value = 26443.666
$ws8.Range("N131").Value = -36523.666

# row 137
$ws8.Range("H137").Value = 45715
$ws8.Range("J137").Value = 45715
$ws8.Range("L137").Value = 45715
$ws8.Range("N137").Value = -55915
